$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.627.66"
$ws.Range("E2").Value = "  -0.85%  "
$ws.Range("D3").Value = "2.226.31"
$ws.Range("E3").Value = "  -0.81%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "252.17"
$ws.Range("E5").Value = "  +8.36%  "
$ws.Range("D6").Value = "0.629"
$ws.Range("E6").Value = "  -1.11%  "
$ws.Range("D7").Value = "71.18"
$ws.Range("E7").Value = "  +1.88%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").Value = "0.569"
$ws.Range("E9").Value = "  +1.83%  "
$ws.Range("D10").Value = "42.61"
$ws.Range("E10").Value = "  +19.17%  "
$ws.Range("D11").Value = "0.0962"
$ws.Range("E11").Value = "  -2.99%  "
$ws.Range("D12").Value = "58.75"
$ws.Range("E12").Value = "  +0.60%  "
$ws.Range("E13").Value = "  +0.36%  "
$ws.Range("D14").Value = "7.05"
$ws.Range("E14").Value = "  +3.31%  "
$ws.Range("D15").Value = "2.558.84"
$ws.Range("E15").Value = "  -0.60%  "
$ws.Range("E16").Value = "  -1.44%  "
$ws.Range("E17").Value = "  -1.11%  "
$ws.Range("D18").Value = "2.226.47"
$ws.Range("E18").Value = "  -0.75%  "
$ws.Range("D19").Value = "41.573.33"
$ws.Range("E19").Value = "  -0.74%  "
$ws.Range("D20").Value = "0.0₃0964"
$ws.Range("E20").Value = "  -1.60%  "
$ws.Range("E21").Value = "  -0.80%  "
$ws.Range("D22").Value = "72.95"
$ws.Range("E22").Value = "  -0.78%  "
$ws.Range("D23").Value = "2.28"
$ws.Range("E23").Value = "  +12.41%  "
$ws.Range("D24").Value = "234.56"
$ws.Range("E24").Value = "  -1.34%  "
$ws.Range("D25").Value = "3.88"
$ws.Range("E25").Value = "  +6.68%  "
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("E27").Value = "  +6.27%  "
$ws.Range("D28").Value = "10.44"
$ws.Range("E28").Value = "  +3.69%  "
$ws.Range("E29").Value = "  +1.40%  "
$ws.Range("D30").Value = "171.63"
$ws.Range("E30").Value = "  +1.32%  "
$ws.Range("D31").Value = "20.59"
$ws.Range("E31").Value = "  -0.43%  "
$ws.Range("E32").Value = "  +1.42%  "
$ws.Range("E33").Value = "  -1.92%  "
$ws.Range("E34").Value = "  +1.59%  "
$ws.Range("D35").Value = "0.0721"
$ws.Range("E35").Value = "  +0.14%  "
$ws.Range("D36").Value = "26.70"
$ws.Range("E36").Value = "  +20.82%  "
$ws.Range("E37").Value = "  -2.37%  "
$ws.Range("E38").Value = "  +11.08%  "
$ws.Range("E39").Value = "  +10.01%  "
$ws.Range("D40").Value = "2.29"
$ws.Range("E40").Value = "  +1.19%  "
$ws.Range("D41").Value = "69.48"
$ws.Range("E41").Value = "  +3.76%  "
$ws.Range("D42").Value = "6.00"
$ws.Range("E42").Value = "  -0.80%  "
$ws.Range("D43").Value = "12.33"
$ws.Range("E43").Value = "  +24.20%  "
$ws.Range("E44").Value = "  +9.69%  "
$ws.Range("E45").Value = "  +3.20%  "
$ws.Range("B46").Value = "SynthetixNetwork"
$ws.Range("C46").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$ws.Range("D46").Value = "4.81"
$ws.Range("E46").Value = "  +10.17%  "
$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").Value = "8.81"
$ws.Range("E47").Value = "  -3.15%  "
$ws.Range("E48").Value = "  +1.12%  "
$ws.Range("E49").Value = "  +0.37%  "
$ws.Range("E50").Value = "  +6.92%  "
$ws.Range("D51").Value = "1.20"
$ws.Range("E51").Value = "  +2.13%  "
